$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '28.242.33'
$ws.Cells.Item(2, 5).Value = '  +1.22%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.805.80'
$ws.Cells.Item(3, 5).Value = '  +3.05%  '

$ws.Cells.Item(4, 5).Value = '  -0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '335.65'
$ws.Cells.Item(5, 5).Value = '  -0.13%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.9996'
$ws.Cells.Item(6, 5).Value = '  -0.12%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4610'
$ws.Cells.Item(7, 5).Value = '  +20.37%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3731'
$ws.Cells.Item(8, 5).Value = '  +9.53%  '

$ws.Cells.Item(9, 5).Value = '  -1.65%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.07646'
$ws.Cells.Item(10, 5).Value = '  +5.74%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '1.151'
$ws.Cells.Item(11, 5).Value = '  +3.14%  '

$ws.Cells.Item(12, 2).Value = 'BinanceUSD'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.002'
$ws.Cells.Item(12, 5).Value = '  -0.05%  '

$ws.Cells.Item(13, 2).Value = 'Solana'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '22.37'
$ws.Cells.Item(13, 5).Value = '  -0.68%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.341'
$ws.Cells.Item(14, 5).Value = '  +2.81%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.474'
$ws.Cells.Item(15, 5).Value = '  +4.61%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.806.54'
$ws.Cells.Item(16, 5).Value = '  +3.01%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.00001098'
$ws.Cells.Item(17, 5).Value = '  +3.41%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.06720'
$ws.Cells.Item(18, 5).Value = '  +1.54%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '81.89'
$ws.Cells.Item(19, 5).Value = '  +3.61%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.9995'
$ws.Cells.Item(20, 5).Value = '  -0.11%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '17.44'
$ws.Cells.Item(21, 5).Value = '  +4.30%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.418'
$ws.Cells.Item(22, 5).Value = '  +3.12%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '28.238.92'
$ws.Cells.Item(23, 5).Value = '  +1.20%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '11.86'
$ws.Cells.Item(24, 5).Value = '  +1.62%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.409'
$ws.Cells.Item(25, 5).Value = '  +1.12%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '20.82'
$ws.Cells.Item(26, 5).Value = '  +4.71%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '153.90'
$ws.Cells.Item(27, 5).Value = '  +0.03%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.378'
$ws.Cells.Item(28, 5).Value = '  +2.85%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.011.59'
$ws.Cells.Item(29, 5).Value = '  +2.95%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '133.35'
$ws.Cells.Item(30, 5).Value = '  +0.97%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.255'
$ws.Cells.Item(31, 5).Value = '  -1.12%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '4.027'
$ws.Cells.Item(32, 5).Value = '  +0.11%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.09568'
$ws.Cells.Item(33, 5).Value = '  +8.38%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '5.866'
$ws.Cells.Item(34, 5).Value = '  +0.31%  '

$ws.Cells.Item(35, 5).Value = '  +5.68%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '12.12'
$ws.Cells.Item(36, 5).Value = '  -0.88%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.06363'
$ws.Cells.Item(37, 5).Value = '  +3.30%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02352'
$ws.Cells.Item(38, 5).Value = '  +2.76%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '5.253'
$ws.Cells.Item(39, 5).Value = '  +1.87%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.6638'
$ws.Cells.Item(40, 5).Value = '  +1.00%  '

$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.237'
$ws.Cells.Item(41, 5).Value = '  +2.05%  '

$ws.Cells.Item(42, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.499'
$ws.Cells.Item(42, 5).Value = '  -0.14%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '8.258'
$ws.Cells.Item(43, 5).Value = '  +3.19%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '14.41'
$ws.Cells.Item(44, 5).Value = '  +5.26%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.9991'
$ws.Cells.Item(45, 5).Value = '  -0.13%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.6115'
$ws.Cells.Item(46, 5).Value = '  +0.75%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '3.824'
$ws.Cells.Item(47, 5).Value = '  -0.09%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '129.64'
$ws.Cells.Item(48, 5).Value = '  +2.39%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.049'
$ws.Cells.Item(49, 5).Value = '  +2.17%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.07164'
$ws.Cells.Item(50, 5).Value = '  +2.70%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.178'
$ws.Cells.Item(51, 5).Value = '  +0.43%  '
